# Applies the cryptos-list refresh described in the commit diff.
# For numeric-looking "Price" strings (e.g. "1.006"), Excel's Value setter
# would auto-convert them to real numbers. We force those specific cells to
# text via NumberFormat "@" and then restore the default "Normal" style so
# only the cell value changes (matching the original inline-string cells).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "27.461.22"
$ws.Range("E2").Value = "  -0.55%  "
$ws.Range("E3").Value = "  -1.95%  "
$ws.Range("D4").NumberFormat = "@"
$ws.Range("D4").Value = "1.006"
$ws.Range("D4").Style = "Normal"
$ws.Range("E4").Value = "  -0.44%  "
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "332.66"
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = "  -0.42%  "
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "1.005"
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = "  -0.56%  "
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = "0.4577"
$ws.Range("D7").Style = "Normal"
$ws.Range("E7").Value = "  -2.30%  "
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "0.3803"
$ws.Range("D8").Style = "Normal"
$ws.Range("E8").Value = "  -2.68%  "
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "46.37"
$ws.Range("D9").Style = "Normal"
$ws.Range("E9").Value = "  +1.51%  "
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "0.07874"
$ws.Range("D10").Style = "Normal"
$ws.Range("E10").Value = "  -1.35%  "
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "0.9688"
$ws.Range("D11").Style = "Normal"
$ws.Range("E11").Value = "  -3.24%  "
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "21.01"
$ws.Range("D12").Style = "Normal"
$ws.Range("E12").Value = "  -3.35%  "
$ws.Range("D13").Value = "1.826.28"
$ws.Range("E13").Value = "  -2.59%  "
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "5.884"
$ws.Range("D14").Style = "Normal"
$ws.Range("E14").Value = "  -1.74%  "
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "7.043"
$ws.Range("D15").Style = "Normal"
$ws.Range("E15").Value = "  -2.66%  "
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "1.006"
$ws.Range("D16").Style = "Normal"
$ws.Range("E16").Value = "  -0.69%  "
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "89.55"
$ws.Range("D17").Style = "Normal"
$ws.Range("E17").Value = "  +1.65%  "
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "0.06643"
$ws.Range("D18").Style = "Normal"
$ws.Range("E18").Value = "  -0.99%  "
$ws.Range("E19").Value = "  -1.69%  "
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "17.08"
$ws.Range("D20").Style = "Normal"
$ws.Range("E20").Value = "  +0.17%  "
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "1.005"
$ws.Range("D21").Style = "Normal"
$ws.Range("E21").Value = "  -0.47%  "
$ws.Range("D22").Value = "27.444.90"
$ws.Range("E22").Value = "  -0.57%  "
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "5.330"
$ws.Range("D23").Style = "Normal"
$ws.Range("E23").Value = "  -2.22%  "
$ws.Range("E24").Value = "  -0.69%  "
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "2.302"
$ws.Range("D25").Style = "Normal"
$ws.Range("E25").Value = "  -0.30%  "
$ws.Range("D26").Value = "2.042.70"
$ws.Range("E26").Value = "  -2.53%  "
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "155.77"
$ws.Range("D27").Style = "Normal"
$ws.Range("E27").Value = "  -2.48%  "
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "19.36"
$ws.Range("D28").Style = "Normal"
$ws.Range("E28").Value = "  -2.21%  "
$ws.Range("E29").Value = "  -3.97%  "
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "5.268"
$ws.Range("D30").Style = "Normal"
$ws.Range("E30").Value = "  -2.76%  "
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "118.21"
$ws.Range("D31").Style = "Normal"
$ws.Range("E31").Value = "  -2.59%  "
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "0.9413"
$ws.Range("D32").Style = "Normal"
$ws.Range("E32").Value = "  -3.39%  "
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "0.09300"
$ws.Range("D33").Style = "Normal"
$ws.Range("E33").Value = "  -1.73%  "
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "3.589"
$ws.Range("D34").Style = "Normal"
$ws.Range("E34").Value = "  -0.76%  "
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "5.245"
$ws.Range("D35").Style = "Normal"
$ws.Range("E35").Value = "  -0.82%  "
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "1.318"
$ws.Range("D36").Style = "Normal"
$ws.Range("E36").Value = "  -0.99%  "
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "0.05913"
$ws.Range("D37").Style = "Normal"
$ws.Range("E37").Value = "  -2.33%  "
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "0.02184"
$ws.Range("D38").Style = "Normal"
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "8.048"
$ws.Range("D39").Style = "Normal"
$ws.Range("E39").Value = "  -3.09%  "
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "1.156"
$ws.Range("D40").Style = "Normal"
$ws.Range("E40").Value = "  -2.96%  "
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "0.5755"
$ws.Range("D41").Style = "Normal"
$ws.Range("E41").Value = "  -3.05%  "
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "0.1826"
$ws.Range("D42").Style = "Normal"
$ws.Range("E42").Value = "  -2.96%  "
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "9.970"
$ws.Range("D43").Style = "Normal"
$ws.Range("E43").Value = "  -2.58%  "
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "1.267"
$ws.Range("D44").Style = "Normal"
$ws.Range("E44").Value = "  +1.12%  "
$ws.Range("B45").Value = "Decentraland"
$ws.Range("C45").Value = "https://coinranking.com/coin/tEf7-dnwV3BXS+decentraland-mana"
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "0.5436"
$ws.Range("D45").Style = "Normal"
$ws.Range("E45").Value = "  -3.42%  "
$ws.Range("B46").Value = "EnergySwap"
$ws.Range("C46").Value = "https://coinranking.com/coin/SbWqqTui-+energyswap-ens"
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "11.91"
$ws.Range("D46").Style = "Normal"
$ws.Range("E46").Value = "  -2.34%  "
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "1.863"
$ws.Range("D47").Style = "Normal"
$ws.Range("E47").Value = "  -2.98%  "
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "110.80"
$ws.Range("D48").Style = "Normal"
$ws.Range("E48").Value = "  -1.19%  "
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "0.06603"
$ws.Range("D49").Style = "Normal"
$ws.Range("E49").Value = "  -2.34%  "
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "1.005"
$ws.Range("D50").Style = "Normal"
$ws.Range("E50").Value = "  -0.63%  "
$ws.Range("E51").Value = "  -1.20%  "
